$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 271, shifting existing rows 271-402 down to 272-403
$ws.Rows.Item(271).Insert()

# Populate the newly inserted row 271 with its data
$ws.Range("A271").Value = 5
$ws.Range("B271").Value = "Macroferia Regional de Talca"
$ws.Range("C271").Value = "Maule"
$ws.Range("D271").Value = 44806
$ws.Range("E271").Value = 7
$ws.Range("F271").Value = 100112023
$ws.Range("G271").Value = "Brócoli"
$ws.Range("H271").Value = "Sin especificar"
$ws.Range("I271").Value = "Primera"
$ws.Range("J271").Value = 3000
$ws.Range("K271").Value = 1200
$ws.Range("L271").Value = 1200
$ws.Range("M271").Value = 1200
$ws.Range("N271").Value = "$/unidad"
$ws.Range("O271").Value = "Región del Maule"
$ws.Range("P271").Value = 1200
$ws.Range("Q271").Value = 1
$ws.Range("R271").Value = "Hortaliza"
